$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1/K1 were shared-string text ("r"/"s"); the retrained model now writes
# numeric weights there, matching the rest of the column.
$ws.Range("J1").Value = 0.6
$ws.Range("K1").Value = 0.3

# Rows 2-51: J was 1 -> 0.6, K was 0.6 -> 0.3 (model retrained with all data)
$ws.Range("J2:J51").Value = 0.6
$ws.Range("K2:K51").Value = 0.3

# View state: scroll to row 23 / column A, zoom to 100%, and select K1:K51
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("K1:K51").Select() | Out-Null
